$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSheet")

# Remove header cells G1:Z1 (shrink the used columns from A:Z to A:F)
$ws.Range("G1:Z1").Clear()

# Update row 3 values: A3,B3,C3 2 -> 3 (D3 stays 3)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3

# Update row 4 values: 4 -> 5
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5

# Update row 5 values: 5 -> 7
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 7

# Add new rows 6-10 continuing the arithmetic progression (step 2)
$values = @(9, 11, 13, 15, 17)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 6 + $i
    $val = $values[$i]
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
    $ws.Cells.Item($row, 3).Value = $val
    $ws.Cells.Item($row, 4).Value = $val
}

# Update the selected cell
$ws.Range("F3").Select()
